# Update the exercise name ("Cập nhật tên bài tập")
#
# 1) Slide 1 (title slide): the subtitle's "Bài " + "24. " runs are
#    retyped as a single "Bài 24. " run (same formatting, so the text
#    engine merges them back into one run when the sub-range is reset).
# 2) Slide 32: the title "Bài tập" becomes "Bài tập 24.1" - only the
#    "tập" run's characters are touched, leaving "Bài" and the space
#    run untouched.

$p = $ppt.ActivePresentation

# --- Slide 1: "Bài " + "24. " -> "Bài 24. " -----------------------------
$slide1 = $p.Slides.Item(1)
$subTitle = $slide1.Shapes.Item(2)
$subTitle.TextFrame.TextRange.Characters(1, 8).Text = "Bài 24. "

# --- Slide 32: "Bài tập" -> "Bài tập 24.1" ------------------------------
$slide32 = $p.Slides.Item(32)
$title32 = $slide32.Shapes.Item(2)
$title32.TextFrame.TextRange.Characters(5, 3).Text = "tập 24.1"
